$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-11 13:14:03"
$wsZhCn.Range("G2").Value = "2016-01-11 13:15:10"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-11 13:14:21"
$wsDeDe.Range("G2").Value = "2016-01-11 13:15:40"
